$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (volume/week number, date range) ---
$ws.Range("A8").Value = "Volume 32   Number  22"
$ws.Range("C9").Value = "Report Covering the Week  5/26/2025  Through  6/1/2025"

# --- Column H width adjustment to match widened values ---
$ws.Columns.Item(8).ColumnWidth = $ws.Columns.Item(5).ColumnWidth

# --- Weekly crime statistics grid updates ---
$ws.Range("D14").NumberFormat = $ws.Range("I16").NumberFormat
$ws.Range("D14").Value = 1
$ws.Range("E14").NumberFormat = $ws.Range("K16").NumberFormat
$ws.Range("E14").Value = -100
$ws.Range("F14").NumberFormat = $ws.Range("I16").NumberFormat
$ws.Range("F14").Value = 1
$ws.Range("G14").NumberFormat = $ws.Range("I16").NumberFormat
$ws.Range("G14").Value = 1
$ws.Range("H14").NumberFormat = $ws.Range("K16").NumberFormat
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 2
$ws.Range("J14").Value = 4
$ws.Range("K14").Value = -50
$ws.Range("L14").Value = -66.666666666666
$ws.Range("M14").Value = 100
$ws.Range("N14").Value = -86.666666666666
$ws.Range("C15").NumberFormat = $ws.Range("A14").NumberFormat
$ws.Range("C15").Value = "0"
$ws.Range("D15").NumberFormat = $ws.Range("A14").NumberFormat
$ws.Range("D15").Value = "0"
$ws.Range("E15").NumberFormat = $ws.Range("A14").NumberFormat
$ws.Range("E15").Value = "***.*"
$ws.Range("F15").Value = 1
$ws.Range("H15").Value = -66.666666666666
$ws.Range("M15").Value = 14.285714285714
$ws.Range("N15").Value = -73.333333333333
$ws.Range("C16").Value = 6
$ws.Range("D16").Value = 8
$ws.Range("E16").Value = -25
$ws.Range("F16").Value = 17
$ws.Range("G16").Value = 25
$ws.Range("H16").Value = -32
$ws.Range("I16").Value = 92
$ws.Range("J16").Value = 99
$ws.Range("K16").Value = -7.070707070707
$ws.Range("L16").Value = 29.577464788732
$ws.Range("M16").Value = -7.070707070707
$ws.Range("N16").Value = -71.25
$ws.Range("D17").Value = 16
$ws.Range("E17").Value = -62.5
$ws.Range("F17").Value = 26
$ws.Range("G17").Value = 35
$ws.Range("H17").Value = -25.714285714285
$ws.Range("I17").Value = 179
$ws.Range("J17").Value = 184
$ws.Range("K17").Value = -2.717391304347
$ws.Range("L17").Value = 20.945945945945
$ws.Range("M17").Value = 138.666666666667
$ws.Range("N17").Value = -21.145374449339
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = -50
$ws.Range("F18").Value = 9
$ws.Range("G18").Value = 9
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 53
$ws.Range("J18").Value = 77
$ws.Range("K18").Value = -31.168831168831
$ws.Range("L18").Value = 1.923076923076
$ws.Range("M18").Value = 130.434782608696
$ws.Range("N18").Value = -67.682926829268
$ws.Range("C19").Value = 5
$ws.Range("D19").Value = 7
$ws.Range("E19").Value = -28.571428571428
$ws.Range("F19").Value = 32
$ws.Range("G19").Value = 37
$ws.Range("H19").Value = -13.513513513513
$ws.Range("I19").Value = 163
$ws.Range("J19").Value = 158
$ws.Range("K19").Value = 3.164556962025
$ws.Range("L19").Value = -5.232558139534
$ws.Range("M19").Value = 71.578947368421
$ws.Range("N19").Value = -34.538152610441
$ws.Range("C20").NumberFormat = $ws.Range("A14").NumberFormat
$ws.Range("C20").Value = "0"
$ws.Range("D20").NumberFormat = $ws.Range("A14").NumberFormat
$ws.Range("D20").Value = "0"
$ws.Range("E20").NumberFormat = $ws.Range("A14").NumberFormat
$ws.Range("E20").Value = "***.*"
$ws.Range("F20").Value = 8
$ws.Range("G20").Value = 3
$ws.Range("H20").Value = 166.666666666667
$ws.Range("N20").Value = -80.503144654088
$ws.Range("C21").Value = 18
$ws.Range("D21").Value = 34
$ws.Range("E21").Value = -47.058823529411
$ws.Range("F21").Value = 94
$ws.Range("H21").Value = -16.814159292035
$ws.Range("I21").Value = 528
$ws.Range("J21").Value = 551
$ws.Range("K21").Value = -4.174228675136
$ws.Range("L21").Value = 6.666666666666
$ws.Range("M21").Value = 66.037735849056
$ws.Range("N21").Value = -54.639175257732
$ws.Range("C23").Value = 7
$ws.Range("D23").Value = 15
$ws.Range("E23").Value = -53.333333333333
$ws.Range("F23").Value = 29
$ws.Range("G23").Value = 36
$ws.Range("H23").Value = -19.444444444444
$ws.Range("I23").Value = 170
$ws.Range("J23").Value = 164
$ws.Range("K23").Value = 3.658536585365
$ws.Range("L23").Value = 9.677419354838
$ws.Range("M23").Value = 84.782608695652
$ws.Range("C24").Value = 8
$ws.Range("D24").Value = 6
$ws.Range("E24").Value = 33.333333333333
$ws.Range("F24").Value = 66
$ws.Range("G24").Value = 49
$ws.Range("H24").Value = 34.693877551020
$ws.Range("I24").Value = 436
$ws.Range("J24").Value = 320
$ws.Range("K24").Value = 36.25
$ws.Range("L24").Value = 18.801089918256
$ws.Range("M24").Value = 67.692307692307
$ws.Range("C25").Value = 1
$ws.Range("D25").Value = 2
$ws.Range("E25").Value = -50
$ws.Range("F25").Value = 10
$ws.Range("G25").Value = 18
$ws.Range("H25").Value = -44.444444444444
$ws.Range("J25").Value = 70
$ws.Range("K25").Value = 94.285714285714
$ws.Range("L25").Value = 43.157894736842
$ws.Range("C26").Value = 6
$ws.Range("D26").Value = 30
$ws.Range("E26").Value = -80
$ws.Range("F26").Value = 42
$ws.Range("G26").Value = 89
$ws.Range("H26").Value = -52.808988764044
$ws.Range("I26").Value = 226
$ws.Range("J26").Value = 330
$ws.Range("K26").Value = -31.515151515151
$ws.Range("L26").Value = 0.892857142857
$ws.Range("M26").Value = -21.527777777777
$ws.Range("C27").NumberFormat = $ws.Range("A14").NumberFormat
$ws.Range("C27").Value = "0"
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = -100
$ws.Range("F27").Value = 1
$ws.Range("G27").Value = 4
$ws.Range("H27").Value = -75
$ws.Range("J27").Value = 10
$ws.Range("K27").Value = -10
$ws.Range("L27").Value = -40
$ws.Range("D28").NumberFormat = $ws.Range("A14").NumberFormat
$ws.Range("D28").Value = "0"
$ws.Range("E28").NumberFormat = $ws.Range("A14").NumberFormat
$ws.Range("E28").Value = "***.*"
$ws.Range("F28").Value = 5
$ws.Range("G28").Value = 4
$ws.Range("H28").Value = 25
$ws.Range("I28").Value = 21
$ws.Range("K28").Value = -8.695652173913
$ws.Range("L28").Value = -4.545454545454
$ws.Range("N29").Value = -88.888888888888
$ws.Range("N30").Value = -87.878787878787
